$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns B and D (data rows) to be stored as text so that
# numeric-looking values like "2.529" keep their literal string
# representation instead of being parsed as numbers.
$ws.Range("B2:B49").NumberFormat = "@"
$ws.Range("D2:D49").NumberFormat = "@"

# Row 1 header: only the "Valor final" -> "Valor Final" label changes.
$ws.Range("D1").Value = 'Valor Final'

$ws.Range("A2").Value = 'Ar condicionado Philco  split inverter  frio 9000 BTU  branco 220V PAC9000IFM15'
$ws.Range("B2").Value = '2.529'
$ws.Range("C2").Value = 24
$ws.Range("D2").Value = '1.897'

$ws.Range("A3").Value = 'Caixa De Som Partybox Stage 320 Jbl - Blpbstage320br'
$ws.Range("B3").Value = '4.501'
$ws.Range("C3").Value = 25
$ws.Range("D3").Value = '3.349'

$ws.Range("A4").Value = 'Console padrão Ps5 Slim Bundle Ratchet & Clank and Returnal Cor Branco'
$ws.Range("B4").Value = '4.399'
$ws.Range("C4").Value = 11
$ws.Range("D4").Value = '3.899'

$ws.Range("A5").Value = 'Cooktop 5 Bocas A Gás Dako Supreme Tripla Chama Preto Bivolt 110/220V'
$ws.Range("B5").Value = '848'
$ws.Range("C5").Value = 22
$ws.Range("D5").Value = '655'

$ws.Range("A6").Value = 'Cooktop Itatiaia Essencial 4 Bocas - Preto'
$ws.Range("B6").Value = '439'
$ws.Range("C6").Value = 12
$ws.Range("D6").Value = '385'

$ws.Range("A7").Value = 'Creatina Monohidratada 500g Soldiers Nutrition Sabor Natural'
$ws.Range("B7").Value = '139'
$ws.Range("C7").Value = 28
$ws.Range("D7").Value = '100'

$ws.Range("A8").Value = 'Creatina Monohidratada 600g 100% Pura Soldiers Nutrition'
$ws.Range("B8").Value = '154'
$ws.Range("C8").Value = 25
$ws.Range("D8").Value = '114'

$ws.Range("A9").Value = 'Fogão  de pé Atlas Atenas Glass  4Q  gás engarrafado 4 queimadores  preto 127V/220V porta com visor 50L'
$ws.Range("B9").Value = '1.039'
$ws.Range("C9").Value = 29
$ws.Range("D9").Value = '733'

$ws.Range("A10").Value = 'Fogão Atlas 4 Bocas Preto Atenas Glass - Bivolt'
$ws.Range("B10").Value = '1.039'
$ws.Range("C10").Value = 29
$ws.Range("D10").Value = '733'

$ws.Range("A11").Value = 'Forno Elétrico Philco Pfe44p Dupla Resistência 44l Cor Preto 110V'
$ws.Range("B11").Value = '439'
$ws.Range("C11").Value = 20
$ws.Range("D11").Value = '349'

$ws.Range("A12").Value = 'Fralda Pampers Pants Ajuste Total Max Xxg 78 Unidades'
$ws.Range("B12").Value = '154'
$ws.Range("C12").Value = 33
$ws.Range("D12").Value = '103'

$ws.Range("A13").Value = 'Impressora Multifuncional 3 Em 1 Ecotank L3250 Preta Epson Cor Preto Bivolt'
$ws.Range("B13").Value = '1.195'
$ws.Range("C13").Value = 9
$ws.Range("D13").Value = '1.079'

$ws.Range("A14").Value = 'Impressora a cor multifuncional Epson EcoTank L3250 com wifi preta 220V'
$ws.Range("B14").Value = '1.299'
$ws.Range("C14").Value = 23
$ws.Range("D14").Value = '989'

$ws.Range("A15").Value = 'Mercado Pago: Point Pro2 - A Maquininha De Cartão + Completa'
$ws.Range("B15").Value = '149'
$ws.Range("C15").Value = 20
$ws.Range("D15").Value = '118'

$ws.Range("A16").Value = 'Micro-ondas Efficient 23 Litros Me23b Branco Electrolux 110v'
$ws.Range("B16").Value = '799'
$ws.Range("C16").Value = 30
$ws.Range("D16").Value = '554'

$ws.Range("A17").Value = 'Micro-ondas Electrolux de bancada Branco com Função Tira Odor e Manter Aquecido 34L MEO44 127v'
$ws.Range("B17").Value = '799'
$ws.Range("C17").Value = 15
$ws.Range("D17").Value = '677'

$ws.Range("A18").Value = 'Micro-ondas Philco Pmo23e 20l Espelhado 1100w Tira Odor 110v 127V'
$ws.Range("B18").Value = '799'
$ws.Range("C18").Value = 24
$ws.Range("D18").Value = '599'

$ws.Range("A19").Value = 'Modulo Taramps Ts400x4 400w 2 Ohms Rca Ts 400x4 4 Canais 100w Amplificador 400rms T400 4 Canais Potencia Taramps Som Para Carro Moto Caminhonete Automotivo'
$ws.Range("B19").Value = '221'
$ws.Range("C19").Value = 15
$ws.Range("D19").Value = '187'

$ws.Range("A20").Value = 'Monitor Gamer Samsung T350 24” FHD, Tela Plana, 75Hz, 5ms, HDMI, FreeSync, Game Mode'
$ws.Range("B20").Value = '729'
$ws.Range("C20").Value = 10
$ws.Range("D20").Value = '653'

$ws.Range("A21").Value = 'Motorola Moto G24 128GB Grafite 8GB RAM'
$ws.Range("B21").Value = '1.099'
$ws.Range("C21").Value = 25
$ws.Range("D21").Value = '819'

$ws.Range("A22").Value = 'Motorola Moto G24 128GB Rosa 8GB RAM'
$ws.Range("B22").Value = '999'
$ws.Range("C22").Value = 25
$ws.Range("D22").Value = '747'

$ws.Range("A23").Value = 'Motorola Moto G24 128GB Verde 8GB RAM'
$ws.Range("B23").Value = '999'
$ws.Range("C23").Value = 25
$ws.Range("D23").Value = '747'

$ws.Range("A24").Value = 'Motorola Moto G54 5G 256 GB Azul 8 GB RAM'
$ws.Range("B24").Value = '1.799'
$ws.Range("C24").Value = 35
$ws.Range("D24").Value = '1.169'

$ws.Range("A25").Value = 'Notebook Acer Asp3 A315-510p-34xc I3 8gb 256gb Ssd 15.6 W11'
$ws.Range("B25").Value = '3.799'
$ws.Range("C25").Value = 38
$ws.Range("D25").Value = '2.338'

$ws.Range("A26").Value = 'Notebook Acer Aspire 5 A515-57-55b8 Intel Core I5 8gb 256gb SSD 15,6'''' W11'
$ws.Range("B26").Value = '4.331'
$ws.Range("C26").Value = 35
$ws.Range("D26").Value = '2.789'

$ws.Range("A27").Value = 'Notebook Lenovo Ideapad Celeron 4gb 128ssd 15.6 W11 C/office Cor Cinza'
$ws.Range("B27").Value = '2.665'
$ws.Range("C27").Value = 32
$ws.Range("D27").Value = '1.799'

$ws.Range("A28").Value = 'Notebook Samsung Galaxy Book2 I5-1235u Windows 11 Home 8gb 256gb Ssd Grafite'
$ws.Range("B28").Value = '4.665'
$ws.Range("C28").Value = 37
$ws.Range("D28").Value = '2.899'

$ws.Range("A29").Value = 'Paco Rabanne One million 1 Million Tradicional EDT 200ml para masculino'
$ws.Range("B29").Value = '848'
$ws.Range("C29").Value = 35
$ws.Range("D29").Value = '548'

$ws.Range("A30").Value = 'Parafusadeira E Furadeira Wap Impacto 21v 3/8 K21 Id02 Nova Cor Amarelo Frequência 50/60 Hz 110V/220V'
$ws.Range("B30").Value = '348'
$ws.Range("C30").Value = 20
$ws.Range("D30").Value = '279'

$ws.Range("A31").Value = 'Parafusadeira Furadeira De Impacto Profissional 21v  Modelo TB-21PX 2 Baterias Com Maleta Tb21px The Black Tools'
$ws.Range("B31").Value = '374'
$ws.Range("C31").Value = 22
$ws.Range("D31").Value = '289'

$ws.Range("A32").Value = 'Parafusadeira Furadeira Sem Fio Bateria 12v P/ Madeira Metal Cor Amarelo/Preto Frequência 60 110V/220V'
$ws.Range("B32").Value = '210'
$ws.Range("C32").Value = 23
$ws.Range("D32").Value = '159'

$ws.Range("A33").Value = 'Philco  PMO23EB Branco 220V'
$ws.Range("B33").Value = '699'
$ws.Range("C33").Value = 26
$ws.Range("D33").Value = '514'

$ws.Range("A34").Value = 'Samsung Galaxy A15 4G Dual SIM 128 GB Azul claro 4 GB RAM'
$ws.Range("B34").Value = '1.415'
$ws.Range("C34").Value = 37
$ws.Range("D34").Value = '889'

$ws.Range("A35").Value = 'Samsung Galaxy A15 4G Dual SIM 128 GB Azul escuro 4 GB RAM'
$ws.Range("B35").Value = '1.415'
$ws.Range("C35").Value = 38
$ws.Range("D35").Value = '869'

$ws.Range("A36").Value = 'Samsung Galaxy A15 Dual SIM 4G 256GB Azul claro 8GB RAM'
$ws.Range("B36").Value = '1.499'
$ws.Range("C36").Value = 28
$ws.Range("D36").Value = '1.079'

$ws.Range("A37").Value = 'Samsung Galaxy A15 Dual SIM 5G 256GB Azul-escuro 8GB RAM'
$ws.Range("B37").Value = '1.799'
$ws.Range("C37").Value = 30
$ws.Range("D37").Value = '1.259'

$ws.Range("A38").Value = 'Samsung Galaxy Tab S9 Fe Wifi, 128gb, 6gb Ram, Tela 10.9 Cor Cinza'
$ws.Range("B38").Value = '2.989'
$ws.Range("C38").Value = 21
$ws.Range("D38").Value = '2.357'

$ws.Range("A39").Value = 'Samsung Smart Tv 43'''' Uhd 4k 43cu7700 2023'
$ws.Range("B39").Value = '3.402'
$ws.Range("C39").Value = 41
$ws.Range("D39").Value = '1.984'

$ws.Range("A40").Value = 'Smart TV LG 32’’ LED HD 32LQ621 Bivolt Preta - Experiência Visual Incrível'
$ws.Range("B40").Value = '1.499'
$ws.Range("C40").Value = 27
$ws.Range("D40").Value = '1.089'

$ws.Range("A41").Value = 'Smart TV LG AI ThinQ 43LM631C0SB LED webOS Full HD 43" 100V/240V'
$ws.Range("B41").Value = '1.865'
$ws.Range("C41").Value = 12
$ws.Range("D41").Value = '1.641'

$ws.Range("A42").Value = 'Smart Tv 43 4k Uhd Thinq Ai 43ur7800 Hdr 10 Pro LG Bivolt'
$ws.Range("B42").Value = '2.998'
$ws.Range("C42").Value = 40
$ws.Range("D42").Value = '1.798'

$ws.Range("A43").Value = 'Smart Tv 43'''' Android Dolby Aws-tv-43-bl-02-a Aiwa Bivolt'
$ws.Range("B43").Value = '1.899'
$ws.Range("C43").Value = 24
$ws.Range("D43").Value = '1.439'

$ws.Range("A44").Value = 'Smart Tv LG 50 Led 4k Uhd Wi-fi Bluetooth  Hdr10 50ur871c0sa Preto'
$ws.Range("B44").Value = '3.498'
$ws.Range("C44").Value = 37
$ws.Range("D44").Value = '2.199'

$ws.Range("A45").Value = 'Smart Tv Led 42'''' Ptv42g6fr2cpf Roku Dolby Audio Preta Philco 110V/220V'
$ws.Range("B45").Value = '1.823'
$ws.Range("C45").Value = 20
$ws.Range("D45").Value = '1.459'

$ws.Range("A46").Value = 'Smartphone Motorola Moto g04s 128GB 8GB Ram Boost Camera 16MP com Moto AI sensor FPS lateral - Grafite'
$ws.Range("B46").Value = '879'
$ws.Range("C46").Value = 20
$ws.Range("D46").Value = '699'

$ws.Range("A47").Value = 'Suplemento em Pó Max Titanium sem Sabor em Pote 3000mg'
$ws.Range("B47").Value = '100'
$ws.Range("C47").Value = 20
$ws.Range("D47").Value = '80'

$ws.Range("A48").Value = 'Tablet Samsung Galaxy Tab A9+ 5g 64gb 4gb Ram Grafite'
$ws.Range("B48").Value = '1.499'
$ws.Range("C48").Value = 16
$ws.Range("D48").Value = '1.249'

$ws.Range("A49").Value = 'Varal De Chão Grande Varal De Roupas 3 Andares Dobrável Cor Azul'
$ws.Range("B49").Value = '120'
$ws.Range("C49").Value = 18
$ws.Range("D49").Value = '97'

# Strip the residual number-format style introduced above so the cell
# styling matches the original (unstyled) data rows.
$ws.Range("B2:B49").ClearFormats()
$ws.Range("D2:D49").ClearFormats()

# Drop the now-unused trailing rows (50-52) so the sheet shrinks to 49 rows.
$ws.Rows("50:52").Delete()
